$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels for existing columns B and C, and add new column D header
$ws.Range("B1").Value = "Person mainly works at home"
$ws.Range("C1").Value = "Person sometimes works at home"
$ws.Range("D1").Value = "Person never works at home"

# Copy style from C1 (existing header) to the new D1 header cell
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Update column B (Person mainly works at home) values
$ws.Range("B2").Value = 10403.09735488775
$ws.Range("B3").Value = 10993.74857470999
$ws.Range("B4").Value = 23656.05514102464
$ws.Range("B5").Value = 26551.27259138576
$ws.Range("B6").Value = 20706.19252224476

# Update column C (Person sometimes works at home) values
$ws.Range("C2").Value = 16392.3596801674
$ws.Range("C3").Value = 17792.16297481974
$ws.Range("C4").Value = 16685.02569941839
$ws.Range("C5").Value = 20882.51114162803
$ws.Range("C6").Value = 24627.92460275977

# New column D (Person never works at home) values
$ws.Range("D2").Value = 168496.1955315974
$ws.Range("D3").Value = 168551.294987407
$ws.Range("D4").Value = 153735.2513517399
$ws.Range("D5").Value = 148296.3044738683
$ws.Range("D6").Value = 155283.733665019
